$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final ordered list of gene symbols for column A (rows 1..130) reflecting
# the timepoint-analysis update: CD8A, CTSB and CD36 removed, and 19 new
# genes (FADS2, GLTP, GSTT1, HAMP, TCL1A, TPPP3, WASL, ATRX, BGN, CDCA8,
# ERBB2, EZR, HNF4A, IGF1, IL16, LMNA, MCM4, NF1, PRMT2) inserted, so the
# final sheet lists 130 gene symbols in column A (net +16 rows).
$genes = @(
    "ACVR1B",
    "B4GALT1",
    "BAX",
    "CCL4",
    "CD8B",
    "COX15",
    "ELOVL6",
    "FSTL1",
    "HIST1H1C",
    "HLA-DOB",
    "IGF1R",
    "JUP",
    "KLF1",
    "MAK",
    "MSC",
    "OLFM1",
    "PGM1",
    "PLA2G16",
    "SMC2",
    "ST14",
    "SYNJ2",
    "TIPIN",
    "TNNT1",
    "TUBB2A",
    "FCER1A",
    "MYOM2",
    "S100B",
    "CD22",
    "COCH",
    "FADS2",
    "FGFR1",
    "GLTP",
    "GSTT1",
    "HAMP",
    "IGLC1",
    "MBP",
    "NKX3-1",
    "PLVAP",
    "QRSL1",
    "RHD",
    "RPL37A",
    "SPON1",
    "TCL1A",
    "TPPP3",
    "VPREB3",
    "WASL",
    "ADIPOR2",
    "ATRX",
    "BCL2A1",
    "BGN",
    "CD40",
    "CDCA8",
    "DDX39A",
    "E2F2",
    "EIF4H",
    "ERBB2",
    "ESR1",
    "EZR",
    "F7",
    "FSTL3",
    "GLIPR1",
    "GRHPR",
    "GYPE",
    "HGF",
    "HNF4A",
    "IDH3B",
    "IGF1",
    "IL16",
    "LMNA",
    "LTB",
    "MARCH6",
    "MCM4",
    "MYH11",
    "NF1",
    "NR3C1",
    "PA2G4",
    "PDLIM5",
    "PRMT2",
    "PRPF31",
    "RPS19",
    "RRAD",
    "SFPQ",
    "SMAD3",
    "SORBS3",
    "TPD52L1",
    "UBE2L3",
    "ATG2A",
    "CAMTA2",
    "DIP2A",
    "HYOU1",
    "OTUD7B",
    "POLD3",
    "PSTPIP1",
    "RRNAD1",
    "SCML1",
    "SMARCD2",
    "TCIRG1",
    "ANK1",
    "CANX",
    "CAPZB",
    "CASP1",
    "CD74",
    "CDK2",
    "CELF2",
    "ELN",
    "EPB42",
    "FBXO9",
    "FLNA",
    "HS2ST1",
    "KYNU",
    "NCBP1",
    "NCOA4",
    "PCLO",
    "PRF1",
    "PRPF4B",
    "RBM5",
    "SLC2A3",
    "SLC4A1",
    "SULT1A1",
    "TGM2",
    "UBE2D1",
    "ZFP36L1",
    "DCAF8",
    "DYNLL1",
    "HLA-B",
    "IL2RG",
    "MUC8",
    "NCOA2",
    "PHLPP1",
    "RUBCNL"
)

for ($i = 0; $i -lt $genes.Length; $i++) {
    $ws.Cells.Item($i + 1, 1).Value = $genes[$i]
}
